$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix capitalization of the class-name labels in column A
$ws.Range("A2").Value = "mdaTextHomePage"
$ws.Range("A4").Value = "mdaTitle"
$ws.Range("A8").Value = "pageTitleNewTab"

# Update the selected cell on the sheet
$ws.Range("A4").Select()
